# BankProjections "metric overrides" sheet: add a counter-item-type column
# (mirrors the existing boolean-flags table with a new "Borrowings" /
# "Agio" column) and rename the "Offset liquidity" item type to
# "Offset pnl".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("metric overrides")

# 1. Insert a new row above the "Date*Metric" header row (row 6) to make
#    room for the "counter item type" / "Borrowings" row. This pushes the
#    header + data rows (and the trailing blank styled rows) down by one.
$ws.Range("A6:E6").EntireRow.Insert()

# 2. Fill in the new row 6: "counter item type" label + "Borrowings" value.
$ws.Range("A6").Value = "counter item type"
$ws.Range("E6").Value = "Borrowings"

# 3. Rename the "Offset liquidity" item type to "Offset pnl".
$ws.Range("A3").Value = "Offset pnl"

# 4. Extend the header row (now row 7) with the new "Agio" column.
$ws.Range("E7").Value = "Agio"

# 5. Populate column E for the boolean flag rows (rows 3-5).
$ws.Range("E3").Value = $false
$ws.Range("E4").Value = $true
$ws.Range("E5").Value = $true

# 6. Extend the three data rows (now rows 8-10) with a value of 0.5 in the
#    new column, matching the percentage format used by column D.
$ws.Range("D8:D10").Copy()
$ws.Range("E8:E10").PasteSpecial(-4122)
$ws.Range("E8").Value = 0.5
$ws.Range("E9").Value = 0.5
$ws.Range("E10").Value = 0.5

# 7. Restore the selection Excel left the sheet in after editing.
$null = $ws.Range("D6").Select()
